$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures per the Feb 17 2023 GitHub Actions refresh.
$updates = @{
    "D2" = "308.89"
    "E2" = "-3.85%"
    "D3" = "54.10"
    "E3" = "9.28%"
    "D4" = "5.108"
    "E4" = "-3.73%"
    "D5" = "0.07821"
    "E5" = "-2.73%"
    "D6" = "4.542"
    "E6" = "-1.23%"
    "D7" = "1.369"
    "E7" = "-2.02%"
    "D8" = "1.626"
    "E8" = "-0.87%"
    "D9" = "0.1232"
    "E9" = "-5.84%"
    "D10" = "0.2024"
    "E10" = "3.61%"
    "D11" = "0.04721"
    "E11" = "0.68%"
    "D12" = "0.09406"
    "E12" = "-0.57%"
    "D13" = "0.1044"
    "E13" = "0.07%"
    "D14" = "0.001260"
    "E14" = "-4.42%"
    "D15" = "0.005803"
    "E15" = "-2.04%"
    "E16" = "2,020.33%"
    "D17" = "3.323"
    "E17" = "-0.54%"
    "D18" = "2.417"
    "E18" = "-2.35%"
    "D19" = "0.3415"
    "D20" = "7.984"
    "E20" = "-1.41%"
    "D21" = "0.1364"
    "E21" = "-0.55%"
    "D22" = "0.2983"
    "E22" = "-3.43%"
    "D23" = "0.04177"
    "E23" = "-0.18%"
    "D24" = "0.001262"
    "E24" = "-3.91%"
    "D25" = "0.003938"
    "E25" = "-8.87%"
    "D26" = "0.0001351"
    "E26" = "0.24%"
    "E38" = "-4.34%"
    "D39" = "0.05877"
    "E39" = "-15.68%"
    "D40" = "0.01101"
    "E40" = "1.49%"
    "D41" = "0.007953"
    "E41" = "-1.18%"
    "D42" = "0.1435"
    "E42" = "-1.47%"
    "D43" = "0.008241"
    "E43" = "5.74%"
    "D44" = "0.008496"
    "E44" = "-1.78%"
    "D45" = "0.3369"
    "E45" = "-3.63%"
    "D46" = "0.00007257"
    "E46" = "9.78%"
    "D47" = "0.00000000751"
    "E47" = "0.32%"
    "D48" = "0.05686"
    "E48" = "-5.84%"
    "D49" = "0.002622"
    "E49" = "-34.30%"
    "D50" = "0.00002102"
    "E50" = "0.32%"
    "D51" = "0.0002002"
    "E51" = "0.32%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"
